$p = $ppt.ActivePresentation

# Content for the new "Title and Content" slides that get appended right
# after the current last slide (146). Slide 155 is added but left blank
# (matching the source deck, where the final new slide has empty
# placeholders).
$newSlides = @(
    @{ Title = "OOP"; Body = "Object-oriented programming has its roots in the 1960s, but it wasn’t until the mid 1980s that it became the main programming paradigm used in the creation of new software. It was developed as a way to handle the rapidly increasing size and complexity of software systems and to make it easier to modify these large and complex systems over time." },
    @{ Title = "CLASS"; Body = "A class in Python is a blueprint for creating objects. Classes encapsulate data for the object and methods to manipulate that data. A class is defined using the class keyword." },
    @{ Title = "Objects"; Body = "An object is an instance of a class. Once a class is defined, you can instantiate it to create an object. Each object can have unique attributes and behavior as defined by its class." },
    @{ Title = "Attributes and Methods"; Body = "Attributes are data stored inside an object or class, and methods are functions defined within a class that operate on its objects. Attributes represent the state of an object, while methods represent the behavior." },
    @{ Title = "Inheritance"; Body = "Inheritance allows a class to inherit attributes and methods from another class, known as the parent class. This promotes code reusability and creates a relationship where the child class can override or extend the functionality of the parent class." },
    @{ Title = "Polymorphism"; Body = "Polymorphism allows for the use of a single interface to represent different underlying forms (data types). In Python, this means that different classes can have methods with the same name, and those methods can be called in the same way even though they might perform different operations." },
    @{ Title = "constructor "; Body = "In object-oriented programming (OOP), a constructor is a special type of subroutine called to create an object. It prepares the new object for use, often accepting arguments that the constructor uses to set required member variables and allocate resources. " },
    @{ Title = "Methods"; Body = "Methods in Python are functions that are defined inside a class and are used to define the behaviors of an object. Unlike standalone functions, methods are called on objects and can access and modify the state of the object to which they belong. This is done through the self parameter, which is a reference to the current instance of the class." },
    @{ Title = ""; Body = "" }
)

$startIndex = $p.Slides.Count + 1

for ($i = 0; $i -lt $newSlides.Count; $i++) {
    $slideIndex = $startIndex + $i
    $s = $p.Slides.Add($slideIndex, 2)

    $item = $newSlides[$i]

    if ($item.Title -ne "") {
        $s.Shapes.Item(1).TextFrame.TextRange.Text = $item.Title
    }

    if ($item.Body -ne "") {
        $s.Shapes.Item(2).TextFrame.TextRange.Text = $item.Body
        $s.Shapes.Item(2).TextFrame.TextRange.Font.Size = 28
    }
}
